$wb = $excel.ActiveWorkbook

# ALC row 32
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 6906.273
$ws.Range("I32").Value = 4650
$ws.Range("J32").Value = 7752.375
$ws.Range("K32").Value = 4650
$ws.Range("L32").Value = 7752.375
$ws.Range("M32").Value = -4324
$ws.Range("N32").Value = -8404.375

# ALC row 51
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 16071.429
$ws.Range("I51").Value = 21500
$ws.Range("J51").Value = 2500
$ws.Range("K51").Value = 21500
$ws.Range("L51").Value = 2500
$ws.Range("M51").Value = -21016
$ws.Range("N51").Value = -3468

# ALC row 70
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 42864140
$ws.Range("I70").Value = 25004638
$ws.Range("K70").Value = 75013914
$ws.Range("M70").Value = -75013644

# ALC row 73
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 42864140
$ws.Range("I73").Value = 25004638
$ws.Range("K73").Value = 75013914
$ws.Range("M73").Value = -75012978

# ALC row 86
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 4500.25
$ws.Range("I86").Value = 3999
$ws.Range("K86").Value = 3999
$ws.Range("M86").Value = -2876

# ALC row 89
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 4500.25
$ws.Range("I89").Value = 3999
$ws.Range("K89").Value = 19995
$ws.Range("M89").Value = -14379

# ALC row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 10636.546
$ws.Range("I113").Value = 14600.5
$ws.Range("J113").Value = 5879.8
$ws.Range("K113").Value = 14600.5
$ws.Range("L113").Value = 5879.8
$ws.Range("M113").Value = -11346.5
$ws.Range("N113").Value = -12387.8

# ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 4970.6
$ws.Range("I116").Value = 5213.25
$ws.Range("J116").Value = 4000
$ws.Range("K116").Value = 5213.25
$ws.Range("L116").Value = 4000
$ws.Range("M116").Value = -1771.25
$ws.Range("N116").Value = -10884

# ALC row 121
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H121").Value = 1736.5
$ws.Range("J121").Value = 1736.5
$ws.Range("L121").Value = 5209.5
$ws.Range("N121").Value = -8703.5

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 8415.929
$ws.Range("I137").Value = 13272.923
$ws.Range("K137").Value = 39818.769
$ws.Range("M137").Value = -37268.769

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1328.9333
$ws.Range("I45").Value = 1020.125
$ws.Range("J45").Value = 1681.8572
$ws.Range("K45").Value = 1020.125
$ws.Range("L45").Value = 1681.8572
$ws.Range("M45").Value = -643.125
$ws.Range("N45").Value = -2435.8572

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4595.852
$ws.Range("I61").Value = 4524.391
$ws.Range("K61").Value = 4524.391
$ws.Range("M61").Value = -4312.391

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1859.1333
$ws.Range("I74").Value = 1917.6923
$ws.Range("J74").Value = 1478.5
$ws.Range("K74").Value = 1917.6923
$ws.Range("L74").Value = 1478.5
$ws.Range("M74").Value = -1043.6923
$ws.Range("N74").Value = -3226.5

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1859.1333
$ws.Range("I77").Value = 1917.6923
$ws.Range("J77").Value = 1478.5
$ws.Range("K77").Value = 9588.461499999999
$ws.Range("L77").Value = 7392.5
$ws.Range("M77").Value = -5220.461499999999
$ws.Range("N77").Value = -16128.5

# ARM row 103
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H103").Value = 50000
$ws.Range("J103").Value = 50000
$ws.Range("L103").Value = 50000
$ws.Range("N103").Value = -52344

# ARM row 109
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H109").Value = 105188.5
$ws.Range("J109").Value = 105188.5
$ws.Range("L109").Value = 105188.5
$ws.Range("N109").Value = -107962.5

# ARM row 117
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H117").Value = 0
$ws.Range("I117").Value = 0
$ws.Range("K117").Value = 0
$ws.Range("M117").ClearContents()

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2637.4285
$ws.Range("I132").Value = 2498.389
$ws.Range("K132").Value = 7495.167
$ws.Range("M132").Value = -4965.167

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 4595.852
$ws.Range("I136").Value = 4524.391
$ws.Range("K136").Value = 13573.173
$ws.Range("M136").Value = -11023.173

# BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2825
$ws.Range("I105").Value = 1866.6666
$ws.Range("K105").Value = 1866.6666
$ws.Range("M105").Value = -119.6666

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2921.9092
$ws.Range("I134").Value = 2877.923
$ws.Range("K134").Value = 8633.769
$ws.Range("M134").Value = -6098.769

# BSM row 140
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H140").Value = 132221.44
$ws.Range("J140").Value = 132221.44
$ws.Range("L140").Value = 132221.44
$ws.Range("N140").Value = -142581.44

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3115.4
$ws.Range("I31").Value = 1688.4117
$ws.Range("J31").Value = 3981.7856
$ws.Range("K31").Value = 1688.4117
$ws.Range("L31").Value = 3981.7856
$ws.Range("M31").Value = -1393.4117
$ws.Range("N31").Value = -4571.7856

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3115.4
$ws.Range("I34").Value = 1688.4117
$ws.Range("J34").Value = 3981.7856
$ws.Range("K34").Value = 1688.4117
$ws.Range("L34").Value = 3981.7856
$ws.Range("M34").Value = -1486.4117
$ws.Range("N34").Value = -4385.7856

# CRP row 86
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 10999
$ws.Range("I86").Value = 10498.5
$ws.Range("K86").Value = 10498.5
$ws.Range("M86").Value = -9375.5

# CRP row 89
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 10999
$ws.Range("I89").Value = 10498.5
$ws.Range("K89").Value = 52492.5
$ws.Range("M89").Value = -46876.5

# CRP row 98
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H98").Value = 37537.93
$ws.Range("J98").Value = 41331.777
$ws.Range("L98").Value = 41331.777
$ws.Range("N98").Value = -45823.777

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2432.0833
$ws.Range("I134").Value = 2198.6365
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 6595.9095
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = -4060.9095
$ws.Range("N134").Value = -20070

# CUL row 3
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 2254.1667
$ws.Range("I3").Value = 2350.0908
$ws.Range("K3").Value = 7050.2724
$ws.Range("M3").Value = -6938.2724

# CUL row 21
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H21").Value = 774.25
$ws.Range("I21").Value = 549.5
$ws.Range("K21").Value = 1648.5
$ws.Range("M21").Value = -1475.5

# CUL row 98
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 327.2
$ws.Range("I98").Value = 312.25
$ws.Range("K98").Value = 936.75
$ws.Range("M98").Value = 561.25

# CUL row 121
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 866.36365
$ws.Range("I121").Value = 552.7143
$ws.Range("K121").Value = 1658.1429
$ws.Range("M121").Value = -348.1428999999998

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2437.8462
$ws.Range("I131").Value = 1220.7778
$ws.Range("K131").Value = 3662.3334
$ws.Range("M131").Value = 1377.6666

# CUL row 137
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 3012.647
$ws.Range("J137").Value = 4493.5
$ws.Range("L137").Value = 13480.5
$ws.Range("N137").Value = -23680.5

# GSM row 75
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H75").Value = 84000
$ws.Range("J75").Value = 84000
$ws.Range("L75").Value = 84000
$ws.Range("N75").Value = -85748

# GSM row 78
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H78").Value = 84000
$ws.Range("J78").Value = 84000
$ws.Range("L78").Value = 252000
$ws.Range("N78").Value = -260736

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I80").Value = 38166.668
$ws.Range("J80").Value = 7739.9
$ws.Range("K80").Value = 38166.668
$ws.Range("L80").Value = 7739.9
$ws.Range("M80").Value = -37168.668
$ws.Range("N80").Value = -9735.9

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I83").Value = 38166.668
$ws.Range("J83").Value = 7739.9
$ws.Range("K83").Value = 190833.34
$ws.Range("L83").Value = 38699.5
$ws.Range("M83").Value = -185841.34
$ws.Range("N83").Value = -48683.5

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3723.4285
$ws.Range("I132").Value = 3509.8462
$ws.Range("K132").Value = 10529.5386
$ws.Range("M132").Value = -7999.5386

# GSM row 135
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("N135").ClearContents()

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2888.1555
$ws.Range("I132").Value = 3453.8845
$ws.Range("J132").Value = 2114
$ws.Range("K132").Value = 10361.6535
$ws.Range("L132").Value = 6342
$ws.Range("M132").Value = -7831.6535
$ws.Range("N132").Value = -11402

# LTW row 135
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H135").Value = 65999
$ws.Range("J135").Value = 65999
$ws.Range("L135").Value = 65999
$ws.Range("N135").Value = -76139

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 6109.3
$ws.Range("J136").Value = 6838.8
$ws.Range("L136").Value = 20516.4
$ws.Range("N136").Value = -25616.4

# WVR row 107
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2731.7
$ws.Range("I107").Value = 2839.625
$ws.Range("K107").Value = 8518.875
$ws.Range("M107").Value = -6598.875

# WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 4659.385
$ws.Range("I126").Value = 4812.2
$ws.Range("J126").Value = 4150
$ws.Range("K126").Value = 14436.6
$ws.Range("L126").Value = 12450
$ws.Range("M126").Value = -11966.6
$ws.Range("N126").Value = -17390
